$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.362.94"
$ws.Range("E2").Value = "  +6.11%  "

# Row 3
$ws.Range("D3").Value = "2.022.03"
$ws.Range("E3").Value = "  +7.31%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.62%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.692"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.73%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.17%  "

# Row 9
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "62.53"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +13.99%  "

# Row 10
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.379"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.75%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0762"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.20%  "

# Row 12
$ws.Range("E12").Value = "  +2.38%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +12.28%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.833"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.53%  "

# Row 15
$ws.Range("D15").Value = "2.314.83"
$ws.Range("E15").Value = "  +7.28%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.68%  "

# Row 17
$ws.Range("D17").Value = "2.031.36"
$ws.Range("E17").Value = "  +7.84%  "

# Row 18
$ws.Range("D18").Value = "37.277.71"
$ws.Range("E18").Value = "  +5.86%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.57%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0857"
$ws.Range("E20").Value = "  +4.74%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.50%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "254.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.42%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.60%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.996"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.39%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.37%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.53%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.02%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.78%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.59%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.129"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.19%  "

# Row 31
$ws.Range("B31").Value = "Gas"
$ws.Range("C31").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +70.71%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.06%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0613"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.72%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.23%  "

# Row 35
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0899"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +25.60%  "

# Row 36
$ws.Range("E36").Value = "  +0.32%  "

# Row 37
$ws.Range("E37").Value = "  +1.05%  "

# Row 38
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.44%  "

# Row 39
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.889"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.29%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.09%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.67%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0226"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.26%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.58%  "

# Row 44
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +21.15%  "

# Row 45
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.90%  "

# Row 46
$ws.Range("D46").Value = "1.369.31"
$ws.Range("E46").Value = "  +3.50%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0847"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.91%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.29%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.73%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +18.72%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.44%  "
